# Update countries & provincias Spain
# Refresh COVID country stats (Pais sheet): new totals for several countries,
# re-ranked rows for countries that swapped position (name + stats move together),
# and the "last updated" timestamp in A1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: refresh the "datos actualizados" timestamp
$ws.Range("A1").Value = 'Datos actualizados a 10 de Agosto de 2020 a las 20:16'

# Row 4: Estados Unidos
$ws.Range("B4").Value = 5218584
$ws.Range("C4").Value = 19140
$ws.Range("D4").Value = 2673846
$ws.Range("E4").Value = 2378909
$ws.Range("G4").Value = 212
$ws.Range("H4").Value = 165829

# Row 5: Brasil
$ws.Range("B5").Value = 3039349
$ws.Range("C5").Value = 3767
$ws.Range("E5").Value = 819620
$ws.Range("G5").Value = 133
$ws.Range("H5").Value = 101269

# Row 6: India
$ws.Range("B6").Value = 2266954
$ws.Range("C6").Value = 52817
$ws.Range("D6").Value = 1580269
$ws.Range("E6").Value = 641333
$ws.Range("G6").Value = 886
$ws.Range("H6").Value = 45352

# Row 15: Reino Unido
$ws.Range("G15").Value = 21
$ws.Range("H15").Value = 46526

# Row 21: Turquia
$ws.Range("B21").Value = 241997
$ws.Range("C21").Value = 1193
$ws.Range("D21").Value = 224970
$ws.Range("E21").Value = 11169
$ws.Range("G21").Value = 14
$ws.Range("H21").Value = 5858

# Row 22: Alemania
$ws.Range("B22").Value = 217563
$ws.Range("C22").Value = 282
$ws.Range("E22").Value = 10403

# Row 27: Canada
$ws.Range("B27").Value = 119723
$ws.Range("C27").Value = 272
$ws.Range("D27").Value = 105986
$ws.Range("E27").Value = 4755
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 8982

# Row 31: Ecuador
$ws.Range("B31").Value = 94701
$ws.Range("C31").Value = 242
$ws.Range("D31").Value = 78608
$ws.Range("E31").Value = 10161
$ws.Range("G31").Value = 10
$ws.Range("H31").Value = 5932

# Row 34: Israel
$ws.Range("B34").Value = 84381
$ws.Range("C34").Value = 1379
$ws.Range("D34").Value = 58986
$ws.Range("E34").Value = 24783
$ws.Range("G34").Value = 12
$ws.Range("H34").Value = 612

# Row 60: Azerbaiyan -> Marruecos
$ws.Range("A60").Value = 'Marruecos'
$ws.Range("B60").Value = 34063
$ws.Range("C60").Value = 826
$ws.Range("D60").Value = 24524
$ws.Range("E60").Value = 9023
$ws.Range("G60").Value = 18
$ws.Range("H60").Value = 516

# Row 61: Marruecos -> Azerbaiyan
$ws.Range("A61").Value = 'Azerbaiyan'
$ws.Range("B61").Value = 33647
$ws.Range("C61").Value = 79
$ws.Range("D61").Value = 30642
$ws.Range("E61").Value = 2513
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 492

# Row 80: Bosnia y Herzegovina
$ws.Range("B80").Value = 14498
$ws.Range("C80").Value = 811
$ws.Range("D80").Value = 8159
$ws.Range("E80").Value = 5914
$ws.Range("G80").Value = 31
$ws.Range("H80").Value = 425

# Row 90: Zambia
$ws.Range("B90").Value = 8210
$ws.Range("C90").Value = 125
$ws.Range("D90").Value = 6802
$ws.Range("E90").Value = 1167
$ws.Range("G90").Value = 6
$ws.Range("H90").Value = 241

# Row 98: Albania -> Libano
$ws.Range("A98").Value = 'Libano'
$ws.Range("B98").Value = 6812
$ws.Range("C98").Value = 295
$ws.Range("D98").Value = 2290
$ws.Range("E98").Value = 4442
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = 80

# Row 99: Mauritania -> Albania
$ws.Range("A99").Value = 'Albania'
$ws.Range("B99").Value = 6536
$ws.Range("C99").Value = 125
$ws.Range("D99").Value = 3379
$ws.Range("E99").Value = 2957
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 200

# Row 100: Libano -> Mauritania
$ws.Range("A100").Value = 'Mauritania'
$ws.Range("B100").Value = 6523
$ws.Range("D100").Value = 5527
$ws.Range("E100").Value = 839
$ws.Range("H100").Value = 157

# Row 104: Republica de Yibuti
$ws.Range("B104").Value = 5347
$ws.Range("C104").Value = 3
$ws.Range("D104").Value = 5120
$ws.Range("E104").Value = 168

# Row 105: Maldivas
$ws.Range("B105").Value = 5157
$ws.Range("C105").Value = 116
$ws.Range("D105").Value = 2835
$ws.Range("E105").Value = 2303

# Row 121: Sri Lanka
$ws.Range("B121").Value = 2870
$ws.Range("C121").Value = 26
$ws.Range("E121").Value = 266

# Row 124: Mali
$ws.Range("B124").Value = 2573
$ws.Range("C124").Value = 6
$ws.Range("D124").Value = 1969
$ws.Range("E124").Value = 479

# Row 144: Georgia -> Siria
$ws.Range("A144").Value = 'Siria'
$ws.Range("B144").Value = 1255
$ws.Range("C144").Value = 67
$ws.Range("D144").Value = 364
$ws.Range("E144").Value = 839
$ws.Range("H144").Value = 52

# Row 145: Republica de Chipre
$ws.Range("B145").Value = 1252
$ws.Range("C145").Value = 10
$ws.Range("E145").Value = 363

# Row 146: Liberia -> Georgia
$ws.Range("A146").Value = 'Georgia'
$ws.Range("B146").Value = 1250
$ws.Range("C146").Value = 25
$ws.Range("D146").Value = 1010
$ws.Range("E146").Value = 223
$ws.Range("H146").Value = 17

# Row 147: Gambia -> Liberia
$ws.Range("A147").Value = 'Liberia'
$ws.Range("B147").Value = 1240
$ws.Range("C147").Value = 3
$ws.Range("D147").Value = 725
$ws.Range("E147").Value = 436
$ws.Range("H147").Value = 79

# Row 148: Burkina Faso -> Gambia
$ws.Range("A148").Value = 'Gambia'
$ws.Range("B148").Value = 1235
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 221
$ws.Range("E148").Value = 991
$ws.Range("H148").Value = 23

# Row 149: Siria -> Burkina Faso
$ws.Range("A149").Value = 'Burkina Faso'
$ws.Range("B149").Value = 1204
$ws.Range("C149").Value = 29
$ws.Range("D149").Value = 984
$ws.Range("E149").Value = 166
$ws.Range("H149").Value = 54

# Row 158: Vietnam
$ws.Range("E158").Value = 433
$ws.Range("G158").Value = 4
$ws.Range("H158").Value = 15

# Row 202: Timor Oriental -> Santa Lucia
$ws.Range("A202").Value = 'Santa Lucia'

# Row 203: Santa Lucia -> Timor Oriental
$ws.Range("A203").Value = 'Timor Oriental'

# Row 213: Montserrat -> Islas Malvinas
$ws.Range("A213").Value = 'Islas Malvinas'
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

# Row 214: Islas Malvinas -> Montserrat
$ws.Range("A214").Value = 'Montserrat'
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
